$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Cells.Item(3, 6).Value = 94
$ws.Cells.Item(5, 6).Value = 9459
$ws.Cells.Item(7, 6).Value = 7848
$ws.Cells.Item(10, 6).Value = 37
$ws.Cells.Item(11, 6).Value = 6747
$ws.Cells.Item(13, 6).Value = 487
$ws.Cells.Item(16, 6).Value = 665
$ws.Cells.Item(22, 6).Value = 114
$ws.Cells.Item(23, 6).Value = 10969
$ws.Cells.Item(25, 6).Value = 59
$ws.Cells.Item(26, 6).Value = 2083
$ws.Cells.Item(27, 6).Value = 2713
$ws.Cells.Item(29, 6).Value = 2439
$ws.Cells.Item(34, 6).Value = 2237
$ws.Cells.Item(36, 6).Value = 1513
$ws.Cells.Item(37, 6).Value = 66
$ws.Cells.Item(38, 6).Value = 39
$ws.Cells.Item(39, 6).Value = 5558
$ws.Cells.Item(41, 6).Value = 1230
$ws.Cells.Item(42, 6).Value = 792
$ws.Cells.Item(43, 6).Value = 145
$ws.Cells.Item(44, 6).Value = 180
$ws.Cells.Item(45, 6).Value = 1093
$ws.Cells.Item(46, 6).Value = 1032
$ws.Cells.Item(47, 6).Value = 1452
$ws.Cells.Item(48, 6).Value = 82
$ws.Cells.Item(49, 6).Value = 1116

$ws = $wb.Worksheets.Item("演出")
$ws.Cells.Item(6, 6).Value = 4
$ws.Cells.Item(18, 6).Value = 10
$ws.Cells.Item(20, 6).Value = 48

$ws = $wb.Worksheets.Item("本地生活")
$ws.Cells.Item(2, 6).Value = 130
$ws.Cells.Item(3, 6).Value = 232

$ws = $wb.Worksheets.Item("全部类型")
$ws.Cells.Item(2, 6).Value = 94
$ws.Cells.Item(4, 6).Value = 9459
$ws.Cells.Item(5, 6).Value = 9459
$ws.Cells.Item(7, 6).Value = 130
$ws.Cells.Item(8, 6).Value = 232
$ws.Cells.Item(10, 6).Value = 4
$ws.Cells.Item(11, 6).Value = 7848
$ws.Cells.Item(14, 6).Value = 37
$ws.Cells.Item(15, 6).Value = 6747
$ws.Cells.Item(16, 6).Value = 6747
$ws.Cells.Item(18, 6).Value = 487
$ws.Cells.Item(20, 6).Value = 665
$ws.Cells.Item(28, 6).Value = 10969
$ws.Cells.Item(30, 6).Value = 59
$ws.Cells.Item(31, 6).Value = 2083
$ws.Cells.Item(32, 6).Value = 2714
$ws.Cells.Item(33, 6).Value = 2439
$ws.Cells.Item(36, 6).Value = 10
$ws.Cells.Item(37, 6).Value = 2237
$ws.Cells.Item(39, 6).Value = 1513
$ws.Cells.Item(40, 6).Value = 39
$ws.Cells.Item(41, 6).Value = 5558
$ws.Cells.Item(42, 6).Value = 48
$ws.Cells.Item(43, 6).Value = 1230
$ws.Cells.Item(44, 6).Value = 792
$ws.Cells.Item(45, 6).Value = 145
$ws.Cells.Item(46, 6).Value = 180
$ws.Cells.Item(47, 6).Value = 1093
$ws.Cells.Item(48, 6).Value = 1032
$ws.Cells.Item(49, 6).Value = 1452
$ws.Cells.Item(50, 6).Value = 82
$ws.Cells.Item(51, 6).Value = 1116
